{"js": "// Replace each old three-digit-by-one-digit multiplication expression\n// with its new value. The document contains 25 short text runs of the\n// exact form \"NNN\u00d7N=\" scattered across table cells; every one of them\n// changes in this edit, and the old values are distinct, so a plain\n// exact-text search-and-replace (scoped to the whole body) is safe and\n// unambiguous.\nconst pairs = [\n  [\"853\u00d78=\", \"711\u00d74=\"],\n  [\"638\u00d77=\", \"461\u00d78=\"],\n  [\"364\u00d78=\", \"495\u00d76=\"],\n  [\"726\u00d73=\", \"493\u00d73=\"],\n  [\"418\u00d76=\", \"284\u00d79=\"],\n  [\"864\u00d74=\", \"999\u00d74=\"],\n  [\"508\u00d75=\", \"746\u00d75=\"],\n  [\"258\u00d75=\", \"364\u00d73=\"],\n  [\"975\u00d75=\", \"609\u00d72=\"],\n  [\"570\u00d75=\", \"120\u00d74=\"],\n  [\"655\u00d79=\", \"577\u00d74=\"],\n  [\"615\u00d75=\", \"113\u00d78=\"],\n  [\"219\u00d78=\", \"767\u00d77=\"],\n  [\"289\u00d76=\", \"166\u00d76=\"],\n  [\"777\u00d78=\", \"470\u00d74=\"],\n  [\"684\u00d79=\", \"630\u00d76=\"],\n  [\"773\u00d78=\", \"511\u00d74=\"],\n  [\"121\u00d79=\", \"154\u00d76=\"],\n  [\"756\u00d78=\", \"816\u00d72=\"],\n  [\"435\u00d72=\", \"517\u00d77=\"],\n  [\"809\u00d72=\", \"327\u00d77=\"],\n  [\"436\u00d74=\", \"937\u00d78=\"],\n  [\"525\u00d72=\", \"402\u00d73=\"],\n  [\"977\u00d76=\", \"552\u00d78=\"],\n  [\"468\u00d74=\", \"679\u00d76=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each old three-digit-by-one-digit multiplication expression\n# with its new value. The document contains 25 short text runs of the\n# exact form \"NNN\u00d7N=\" scattered across table cells; every one of them\n# changes in this edit, and the old values are distinct, so exact-text\n# Find/Replace across the whole document body is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"853\u00d78=\", \"711\u00d74=\"),\n    @(\"638\u00d77=\", \"461\u00d78=\"),\n    @(\"364\u00d78=\", \"495\u00d76=\"),\n    @(\"726\u00d73=\", \"493\u00d73=\"),\n    @(\"418\u00d76=\", \"284\u00d79=\"),\n    @(\"864\u00d74=\", \"999\u00d74=\"),\n    @(\"508\u00d75=\", \"746\u00d75=\"),\n    @(\"258\u00d75=\", \"364\u00d73=\"),\n    @(\"975\u00d75=\", \"609\u00d72=\"),\n    @(\"570\u00d75=\", \"120\u00d74=\"),\n    @(\"655\u00d79=\", \"577\u00d74=\"),\n    @(\"615\u00d75=\", \"113\u00d78=\"),\n    @(\"219\u00d78=\", \"767\u00d77=\"),\n    @(\"289\u00d76=\", \"166\u00d76=\"),\n    @(\"777\u00d78=\", \"470\u00d74=\"),\n    @(\"684\u00d79=\", \"630\u00d76=\"),\n    @(\"773\u00d78=\", \"511\u00d74=\"),\n    @(\"121\u00d79=\", \"154\u00d76=\"),\n    @(\"756\u00d78=\", \"816\u00d72=\"),\n    @(\"435\u00d72=\", \"517\u00d77=\"),\n    @(\"809\u00d72=\", \"327\u00d77=\"),\n    @(\"436\u00d74=\", \"937\u00d78=\"),\n    @(\"525\u00d72=\", \"402\u00d73=\"),\n    @(\"977\u00d76=\", \"552\u00d78=\"),\n    @(\"468\u00d74=\", \"679\u00d76=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
